# Quarterly balance-sheet roll-forward:
#  - drop the oldest reporting quarter (column D) and shift everything left
#  - append the newest quarter's figures in the new rightmost column (M)
# Deleting column D lets Excel itself re-flow the per-column widths and
# prune the now-unreferenced shared strings (old period label / old
# publish-date label), exactly mirroring what a human editor does when
# they right-click "Delete" on the obsolete quarter's column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("D").Delete()

# --- new quarter header (row 8: period label, row 9: publish date) -------
# These look like dates ("1402-02-10"), and a plain .Value assignment lets
# Excel's smart-typing coerce an ISO-looking string into a date serial. Route
# it through a text formula, then Copy + PasteSpecial(values) to flatten it
# back to a literal string cell — that keeps the shared-string literal
# while leaving the existing cell style untouched.
$ws.Range("M8").Formula = "=""فصل چهارم منتهی به 1401/12"""
$ws.Range("M8").Copy()
$ws.Range("M8").PasteSpecial(-4163)

$ws.Range("M9").Formula = "=""1402-02-10"""
$ws.Range("M9").Copy()
$ws.Range("M9").PasteSpecial(-4163)

# The "published (revision n)" label one column over is not pure ISO-8601
# (it carries the trailing " (7)"), so Excel leaves it as text on a normal
# .Value assignment.
$ws.Range("I9").Value = "1402-02-10 (7)"

# --- newest quarter's financial figures (column M) ------------------------
$ws.Range("M12").Value = 2370001
$ws.Range("M13").Value = 23222503
$ws.Range("M14").Value = 32045778
$ws.Range("M15").Value = 2250423
$ws.Range("M16").Value = 472859
$ws.Range("M17").Value = 0
$ws.Range("M18").Value = 60361564
$ws.Range("M19").Value = 0
$ws.Range("M20").Value = 10703626
$ws.Range("M21").Value = 0
$ws.Range("M22").Value = 3365149
$ws.Range("M23").Value = 13091
$ws.Range("M24").Value = 40
$ws.Range("M25").Value = 12039
$ws.Range("M26").Value = 14093905
$ws.Range("M27").Value = 74455469
$ws.Range("M29").Value = 759304
$ws.Range("M30").Value = 40
$ws.Range("M31").Value = 2064077
$ws.Range("M32").Value = 6103818
$ws.Range("M33").Value = 110095
$ws.Range("M34").Value = 0
$ws.Range("M35").Value = 0
$ws.Range("M36").Value = 0
$ws.Range("M37").Value = 9037294
$ws.Range("M38").Value = 0
$ws.Range("M39").Value = 40
$ws.Range("M40").Value = 0
$ws.Range("M41").Value = 1397869
$ws.Range("M42").Value = 1397869
$ws.Range("M43").Value = 10435163
$ws.Range("M45").Value = 1789912
$ws.Range("M46").Value = 0
$ws.Range("M47").Value = 0
$ws.Range("M48").Value = -82828
$ws.Range("M49").Value = 19590
$ws.Range("M50").Value = 178991
$ws.Range("M51").Value = 0
$ws.Range("M52").Value = 40
$ws.Range("M53").Value = 0
$ws.Range("M54").Value = 40
$ws.Range("M55").Value = 0
$ws.Range("M56").Value = 62114641
$ws.Range("M57").Value = 64020306
$ws.Range("M58").Value = 74455469

# New rightmost column keeps the same width rhythm as the old one it replaces.
$ws.Columns("M").ColumnWidth = 31
